$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.496.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.500.21'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.18'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '194.66'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.91%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.643'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.95'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000296'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.41'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.057.85'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '594.19'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.691.08'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.94'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.61'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.83%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.498.76'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.980'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.99'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.57'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.79%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.10'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.68'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.46'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.89'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.24'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +7.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.27'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.113'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.10'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.732.77'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.70%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.77%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.62'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.02'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '489.80'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.06%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.15%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.71%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.16%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.79'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.80%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.37'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000242'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.74%  '
